$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New product listing data (TH1MKT) — replaces the previous 12-row table
# with a new 9-row table. Columns: A=code, B=description, C=market,
# D=group, E=sub, F=unit/pack info.
$rows = @(
    @("20064426", "TWISTKO JGNG BKR 70G", "TH1MKT", "1", "1", "RT,(E-1B)"),
    @("20137884", "CHUPA JELLY FUN MIX ", "TH1MKT", "1", "2", "RT,(E-1B)"),
    @("10003922", "ABC KP&GULA+SS.10X30", "TH1MKT", "2", "1", "RT,(E-4B)"),
    @("20138893", "AMO DRMY STRAW 180ML", "TH1MKT", "2", "2", "RT,(E-4B)"),
    @("20025825", "PRONAS KORNETKU 200G", "TH1MKT", "3", "1", "RT,(E-4B)"),
    @("10002350", "ABC KCP MANIS TGG275", "TH1MKT", "3", "2", "RT,(E-2B)"),
    @("10037208", "CLOSE UP MT.FRSH 110", "TH1MKT", "4", "1", "PT,(E-3B)"),
    @("10037405", "C/LANG KAYU PUTIH 30", "TH1MKT", "4", "2", "RT,(E-6B)"),
    @("20114432", "ANTANGIN HTBTSDA 5'S", "TH1MKT", "4", "3", "RT,(E-4B)")
)

# Columns whose values are purely-numeric-looking strings (product codes /
# group numbers) must be written as Text so they keep the shared-string
# "t=s" representation instead of turning into numeric cells.
$textCols = @(1, 4, 5)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $rowData = $rows[$i]
    for ($col = 1; $col -le 6; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        $val = $rowData[$col - 1]
        if ($textCols -contains $col) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $val
    }
}

# The old table had 12 data rows (rows 2-13); the new table only needs 9
# (rows 2-10), so remove the trailing rows.
$ws.Rows("11:13").Delete()

# Column F narrows slightly now that the longest entry is shorter. The
# saved <col width> is the Excel ColumnWidth plus a constant ~5/6 padding,
# so back that off to land on an exact width="11" in the XML.
$ws.Columns("F").ColumnWidth = 11 - 5/6
